$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 48 (ALC)
$ws.Range("H48").Value = 2000
$ws.Range("I48").Value = 2000
$ws.Range("J48").Value = 0
$ws.Range("K48").Value = 6000
$ws.Range("L48").Value = 0
$ws.Range("M48").Value = -5708
$ws.Range("N48").ClearContents()

# Row 56 (ALC)
$ws.Range("H56").Value = 2000
$ws.Range("I56").Value = 2000
$ws.Range("J56").Value = 0
$ws.Range("K56").Value = 6000
$ws.Range("L56").Value = 0
$ws.Range("M56").Value = -5466
$ws.Range("N56").ClearContents()

# Row 113 (ALC)
$ws.Range("H113").Value = 2445.0557
$ws.Range("I113").Value = 2285.4614
$ws.Range("J113").Value = 2860
$ws.Range("K113").Value = 2285.4614
$ws.Range("L113").Value = 2860
$ws.Range("M113").Value = 968.5385999999999
$ws.Range("N113").Value = -9368

# Row 125 (ALC)
$ws.Range("H125").Value = 2156.4443
$ws.Range("I125").Value = 1850
$ws.Range("J125").Value = 2244
$ws.Range("K125").Value = 16650
$ws.Range("L125").Value = 20196
$ws.Range("M125").Value = -14190
$ws.Range("N125").Value = -25116

# Row 137 (ALC)
$ws.Range("H137").Value = 3610.3667
$ws.Range("I137").Value = 3332.3044
$ws.Range("J137").Value = 4524
$ws.Range("K137").Value = 9996.913199999999
$ws.Range("L137").Value = 13572
$ws.Range("M137").Value = -7446.913199999999
$ws.Range("N137").Value = -18672

$ws = $wb.Worksheets.Item("ARM")
# Row 61 (ARM)
$ws.Range("H61").Value = 2812.0715
$ws.Range("I61").Value = 2050
$ws.Range("J61").Value = 3828.1667
$ws.Range("K61").Value = 2050
$ws.Range("L61").Value = 3828.1667
$ws.Range("M61").Value = -1838
$ws.Range("N61").Value = -4252.1667

# Row 136 (ARM)
$ws.Range("H136").Value = 2812.0715
$ws.Range("I136").Value = 2050
$ws.Range("J136").Value = 3828.1667
$ws.Range("K136").Value = 6150
$ws.Range("L136").Value = 11484.5001
$ws.Range("M136").Value = -3600
$ws.Range("N136").Value = -16584.5001

$ws = $wb.Worksheets.Item("BSM")
# Row 107 (BSM)
$ws.Range("H107").Value = 1500.5454
$ws.Range("I107").Value = 1045
$ws.Range("J107").Value = 1760.8572
$ws.Range("K107").Value = 1045
$ws.Range("L107").Value = 1760.8572
$ws.Range("M107").Value = 875
$ws.Range("N107").Value = -5600.8572

$ws = $wb.Worksheets.Item("CRP")
# Row 29 (CRP)
$ws.Range("H29").Value = 3949.5
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 3949.5
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 3949.5
$ws.Range("N29").Value = -4535.5

# Row 134 (CRP)
$ws.Range("H134").Value = 2108.2727
$ws.Range("I134").Value = 1910.1111
$ws.Range("J134").Value = 3000
$ws.Range("K134").Value = 5730.3333
$ws.Range("L134").Value = 9000
$ws.Range("M134").Value = -3195.3333
$ws.Range("N134").Value = -14070

$ws = $wb.Worksheets.Item("CUL")
# Row 5 (CUL)
$ws.Range("H5").Value = 1032.75
$ws.Range("I5").Value = 758.6
$ws.Range("J5").Value = 1489.6666
$ws.Range("K5").Value = 2275.8
$ws.Range("L5").Value = 4468.9998
$ws.Range("M5").Value = -2163.8
$ws.Range("N5").Value = -4692.9998

# Row 23 (CUL)
$ws.Range("H23").Value = 76923180
$ws.Range("I23").Value = 105.57143
$ws.Range("J23").Value = 166666780
$ws.Range("K23").Value = 316.71429
$ws.Range("L23").Value = 500000340
$ws.Range("M23").Value = -81.71429000000001
$ws.Range("N23").Value = -500000810

# Row 109 (CUL)
$ws.Range("H109").Value = 4531.4
$ws.Range("I109").Value = 963.5
$ws.Range("J109").Value = 6910
$ws.Range("K109").Value = 2890.5
$ws.Range("L109").Value = 20730
$ws.Range("M109").Value = -1850.5
$ws.Range("N109").Value = -22810

# Row 114 (CUL)
$ws.Range("H114").Value = 611.5
$ws.Range("I114").Value = 195.28572
$ws.Range("J114").Value = 935.2222
$ws.Range("K114").Value = 585.85716
$ws.Range("L114").Value = 2805.6666
$ws.Range("M114").Value = 2668.14284
$ws.Range("N114").Value = -9313.6666

# Row 117 (CUL)
$ws.Range("H117").Value = 1463.2727
$ws.Range("I117").Value = 1750
$ws.Range("J117").Value = 1399.5555
$ws.Range("K117").Value = 5250
$ws.Range("L117").Value = 4198.666499999999
$ws.Range("M117").Value = -1808
$ws.Range("N117").Value = -11082.6665

# Row 120 (CUL)
$ws.Range("H120").Value = 6000
$ws.Range("I120").Value = 3000
$ws.Range("J120").Value = 12000
$ws.Range("K120").Value = 9000
$ws.Range("L120").Value = 36000
$ws.Range("M120").Value = -4162
$ws.Range("N120").Value = -45676

# Row 122 (CUL)
$ws.Range("H122").Value = 4731.5415
$ws.Range("I122").Value = 347.4375
$ws.Range("J122").Value = 13499.75
$ws.Range("K122").Value = 3126.9375
$ws.Range("L122").Value = 121497.75
$ws.Range("M122").Value = -676.9375
$ws.Range("N122").Value = -126397.75

# Row 135 (CUL)
$ws.Range("H135").Value = 1032.75
$ws.Range("I135").Value = 758.6
$ws.Range("J135").Value = 1489.6666
$ws.Range("K135").Value = 6827.400000000001
$ws.Range("L135").Value = 13406.9994
$ws.Range("M135").Value = -4292.400000000001
$ws.Range("N135").Value = -18476.9994

# Row 136 (CUL)
$ws.Range("H136").Value = 2455
$ws.Range("I136").Value = 1471.4286
$ws.Range("J136").Value = 2984.6155
$ws.Range("K136").Value = 4414.2858
$ws.Range("L136").Value = 8953.8465
$ws.Range("M136").Value = 685.7142000000003
$ws.Range("N136").Value = -19153.8465

# Row 137 (CUL)
$ws.Range("H137").Value = 18532344
$ws.Range("I137").Value = 20848512
$ws.Range("J137").Value = 3000
$ws.Range("K137").Value = 62545536
$ws.Range("L137").Value = 9000
$ws.Range("M137").Value = -62540436
$ws.Range("N137").Value = -19200

# Row 140 (CUL)
$ws.Range("H140").Value = 1948.96
$ws.Range("I140").Value = 1190.2354
$ws.Range("J140").Value = 3561.25
$ws.Range("K140").Value = 3570.7062
$ws.Range("L140").Value = 10683.75
$ws.Range("M140").Value = 1609.2938
$ws.Range("N140").Value = -21043.75

$ws = $wb.Worksheets.Item("GSM")
# Row 2 (GSM)
$ws.Range("H2").Value = 96
$ws.Range("I2").Value = 80
$ws.Range("J2").Value = 100
$ws.Range("K2").Value = 80
$ws.Range("L2").Value = 100
$ws.Range("M2").Value = 33
$ws.Range("N2").Value = -326

# Row 22 (GSM)
$ws.Range("H22").Value = 63337.668
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 63337.668
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 63337.668
$ws.Range("N22").Value = -64395.668

# Row 43 (GSM)
$ws.Range("H43").Value = 10979.5
$ws.Range("I43").Value = 1960
$ws.Range("J43").Value = 19999
$ws.Range("K43").Value = 1960
$ws.Range("L43").Value = 19999
$ws.Range("M43").Value = -1809
$ws.Range("N43").Value = -20301

# Row 57 (GSM)
$ws.Range("H57").Value = 39686.332
$ws.Range("I57").Value = 0
$ws.Range("J57").Value = 39686.332
$ws.Range("K57").Value = 0
$ws.Range("L57").Value = 39686.332
$ws.Range("N57").Value = -41326.332

# Row 80 (GSM)
$ws.Range("H80").Value = 450855.47
$ws.Range("I80").Value = 753338.5600000001
$ws.Range("J80").Value = 47544.668
$ws.Range("K80").Value = 753338.5600000001
$ws.Range("L80").Value = 47544.668
$ws.Range("M80").Value = -752340.5600000001
$ws.Range("N80").Value = -49540.668

# Row 83 (GSM)
$ws.Range("H83").Value = 450855.47
$ws.Range("I83").Value = 753338.5600000001
$ws.Range("J83").Value = 47544.668
$ws.Range("K83").Value = 3766692.8
$ws.Range("L83").Value = 237723.34
$ws.Range("M83").Value = -3761700.8
$ws.Range("N83").Value = -247707.34

$ws = $wb.Worksheets.Item("LTW")
# Row 7 (LTW)
$ws.Range("H7").Value = 2857.4614
$ws.Range("I7").Value = 1619.5714
$ws.Range("J7").Value = 4301.6665
$ws.Range("K7").Value = 1619.5714
$ws.Range("L7").Value = 4301.6665
$ws.Range("M7").Value = -1507.5714
$ws.Range("N7").Value = -4525.6665

# Row 18 (LTW)
$ws.Range("H18").Value = 58388.332
$ws.Range("I18").Value = 300
$ws.Range("J18").Value = 70006
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 70006
$ws.Range("M18").Value = -128
$ws.Range("N18").Value = -70350

# Row 126 (LTW)
$ws.Range("H126").Value = 2857.4614
$ws.Range("I126").Value = 1619.5714
$ws.Range("J126").Value = 4301.6665
$ws.Range("K126").Value = 4858.7142
$ws.Range("L126").Value = 12904.9995
$ws.Range("M126").Value = -2388.7142
$ws.Range("N126").Value = -17844.9995

$ws = $wb.Worksheets.Item("WVR")
# Row 29 (WVR)
$ws.Range("H29").Value = 94309.38
$ws.Range("I29").Value = 8600
$ws.Range("J29").Value = 380007.34
$ws.Range("K29").Value = 8600
$ws.Range("L29").Value = 380007.34
$ws.Range("M29").Value = -8310
$ws.Range("N29").Value = -380587.34

# Row 81 (WVR)
$ws.Range("H81").Value = 8573.799999999999
$ws.Range("I81").Value = 8717.5
$ws.Range("J81").Value = 7999
$ws.Range("K81").Value = 17435
$ws.Range("L81").Value = 15998
$ws.Range("M81").Value = -16374
$ws.Range("N81").Value = -18120

# Row 84 (WVR)
$ws.Range("H84").Value = 8573.799999999999
$ws.Range("I84").Value = 8717.5
$ws.Range("J84").Value = 7999
$ws.Range("K84").Value = 87175
$ws.Range("L84").Value = 79990
$ws.Range("M84").Value = -81871
$ws.Range("N84").Value = -90598

Write-Host "Applied all edits"